$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2941
$ws.Range("I62").Value = 677.5
$ws.Range("J62").Value = 4450
$ws.Range("K62").Value = 677.5
$ws.Range("L62").Value = 4450
$ws.Range("M62").Value = -53.5
$ws.Range("N62").Value = -5698

$ws.Range("H64").Value = 3351.7856
$ws.Range("I64").Value = 2927.1428
$ws.Range("J64").Value = 3493.3333
$ws.Range("K64").Value = 2927.1428
$ws.Range("L64").Value = 3493.3333
$ws.Range("M64").Value = -2679.1428
$ws.Range("N64").Value = -3989.3333

$ws.Range("H65").Value = 2941
$ws.Range("I65").Value = 677.5
$ws.Range("J65").Value = 4450
$ws.Range("K65").Value = 3387.5
$ws.Range("L65").Value = 22250
$ws.Range("M65").Value = -267.5
$ws.Range("N65").Value = -28490

$ws.Range("H67").Value = 3351.7856
$ws.Range("I67").Value = 2927.1428
$ws.Range("J67").Value = 3493.3333
$ws.Range("K67").Value = 2927.1428
$ws.Range("L67").Value = 3493.3333
$ws.Range("M67").Value = -2069.1428
$ws.Range("N67").Value = -5209.3333

$ws.Range("H116").Value = 2726.92
$ws.Range("I116").Value = 2191
$ws.Range("J116").Value = 3409
$ws.Range("K116").Value = 2191
$ws.Range("L116").Value = 3409
$ws.Range("M116").Value = 1251
$ws.Range("N116").Value = -10293

$ws.Range("H137").Value = 1264.1852
$ws.Range("I137").Value = 1203.125
$ws.Range("J137").Value = 1353
$ws.Range("K137").Value = 3609.375
$ws.Range("L137").Value = 4059
$ws.Range("M137").Value = -1059.375
$ws.Range("N137").Value = -9159

$ws.Range("H138").Value = 2184.7395
$ws.Range("I138").Value = 1423
$ws.Range("J138").Value = 2293.5596
$ws.Range("K138").Value = 4269
$ws.Range("L138").Value = 6880.6788
$ws.Range("M138").Value = 871
$ws.Range("N138").Value = -17160.6788

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 715.9394
$ws.Range("I74").Value = 690.25806
$ws.Range("J74").Value = 1114
$ws.Range("K74").Value = 690.25806
$ws.Range("L74").Value = 1114
$ws.Range("M74").Value = 183.74194
$ws.Range("N74").Value = -2862

$ws.Range("H77").Value = 715.9394
$ws.Range("I77").Value = 690.25806
$ws.Range("J77").Value = 1114
$ws.Range("K77").Value = 3451.2903
$ws.Range("L77").Value = 5570
$ws.Range("M77").Value = 916.7096999999999
$ws.Range("N77").Value = -14306

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6890.75
$ws.Range("I134").Value = 1826
$ws.Range("J134").Value = 27149.75
$ws.Range("K134").Value = 5478
$ws.Range("L134").Value = 81449.25
$ws.Range("M134").Value = -2943
$ws.Range("N134").Value = -86519.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 869.5762999999999
$ws.Range("I31").Value = 739.8409
$ws.Range("J31").Value = 1250.1333
$ws.Range("K31").Value = 739.8409
$ws.Range("L31").Value = 1250.1333
$ws.Range("M31").Value = -444.8409
$ws.Range("N31").Value = -1840.1333

$ws.Range("H34").Value = 869.5762999999999
$ws.Range("I34").Value = 739.8409
$ws.Range("J34").Value = 1250.1333
$ws.Range("K34").Value = 739.8409
$ws.Range("L34").Value = 1250.1333
$ws.Range("M34").Value = -537.8409
$ws.Range("N34").Value = -1654.1333

$ws.Range("H58").Value = 953.1053000000001
$ws.Range("I58").Value = 940.6
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 940.6
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -737.6
$ws.Range("N58").Value = -1406

$ws.Range("H132").Value = 4880.0312
$ws.Range("I132").Value = 5912.7617
$ws.Range("J132").Value = 2908.4546
$ws.Range("K132").Value = 17738.2851
$ws.Range("L132").Value = 8725.363799999999
$ws.Range("M132").Value = -15208.2851
$ws.Range("N132").Value = -13785.3638

$ws.Range("H133").Value = 63285.832
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 63285.832
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 63285.832
$ws.Range("N133").Value = -68345.83199999999

$ws.Range("H134").Value = 7753050
$ws.Range("I134").Value = 8548005
$ws.Range("J134").Value = 2239.5
$ws.Range("K134").Value = 25644015
$ws.Range("L134").Value = 6718.5
$ws.Range("M134").Value = -25641480
$ws.Range("N134").Value = -11788.5

$ws.Range("H136").Value = 953.1053000000001
$ws.Range("I136").Value = 940.6
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2821.8
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -271.8000000000002
$ws.Range("N136").Value = -8100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1999
$ws.Range("I62").Value = 1999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5997
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5311
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 2498.5
$ws.Range("I63").Value = 1001.3333
$ws.Range("J63").Value = 6990
$ws.Range("K63").Value = 3003.9999
$ws.Range("L63").Value = 20970
$ws.Range("M63").Value = -2254.9999
$ws.Range("N63").Value = -22468

$ws.Range("H64").Value = 3839.0908
$ws.Range("I64").Value = 1667.3334
$ws.Range("J64").Value = 4653.5
$ws.Range("K64").Value = 5002.0002
$ws.Range("L64").Value = 13960.5
$ws.Range("M64").Value = -4732.0002
$ws.Range("N64").Value = -14500.5

$ws.Range("H65").Value = 1999
$ws.Range("I65").Value = 1999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17991
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -14559
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 2498.5
$ws.Range("I66").Value = 1001.3333
$ws.Range("J66").Value = 6990
$ws.Range("K66").Value = 9011.9997
$ws.Range("L66").Value = 62910
$ws.Range("M66").Value = -5267.9997
$ws.Range("N66").Value = -70398

$ws.Range("H67").Value = 3839.0908
$ws.Range("I67").Value = 1667.3334
$ws.Range("J67").Value = 4653.5
$ws.Range("K67").Value = 5002.0002
$ws.Range("L67").Value = 13960.5
$ws.Range("M67").Value = -4066.0002
$ws.Range("N67").Value = -15832.5

$ws.Range("H68").Value = 1362.7354
$ws.Range("I68").Value = 699.61536
$ws.Range("J68").Value = 1773.238
$ws.Range("K68").Value = 2098.84608
$ws.Range("L68").Value = 5319.714
$ws.Range("M68").Value = -1287.84608
$ws.Range("N68").Value = -6941.714

$ws.Range("H71").Value = 1362.7354
$ws.Range("I71").Value = 699.61536
$ws.Range("J71").Value = 1773.238
$ws.Range("K71").Value = 6296.53824
$ws.Range("L71").Value = 15959.142
$ws.Range("M71").Value = -2240.53824
$ws.Range("N71").Value = -24071.142

$ws.Range("H98").Value = 771.3333
$ws.Range("I98").Value = 112.4
$ws.Range("J98").Value = 1595
$ws.Range("K98").Value = 337.2
$ws.Range("L98").Value = 4785
$ws.Range("M98").Value = 1160.8
$ws.Range("N98").Value = -7781

$ws.Range("H107").Value = 4751.04
$ws.Range("I107").Value = 620.5833
$ws.Range("J107").Value = 8563.77
$ws.Range("K107").Value = 1861.7499
$ws.Range("L107").Value = 25691.31
$ws.Range("M107").Value = 58.25009999999997
$ws.Range("N107").Value = -29531.31

$ws.Range("H131").Value = 18519826
$ws.Range("I131").Value = 142857870
$ws.Range("J131").Value = 1394.5106
$ws.Range("K131").Value = 428573610
$ws.Range("L131").Value = 4183.531800000001
$ws.Range("M131").Value = -428568570
$ws.Range("N131").Value = -14263.5318

$ws.Range("H133").Value = 3288.7407
$ws.Range("I133").Value = 1477
$ws.Range("J133").Value = 3515.2083
$ws.Range("K133").Value = 4431
$ws.Range("L133").Value = 10545.6249
$ws.Range("M133").Value = 629
$ws.Range("N133").Value = -20665.6249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 36739.4
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 36739.4
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 36739.4
$ws.Range("N133").Value = -46859.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2026.3
$ws.Range("I7").Value = 1917.5555
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 1917.5555
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -1805.5555
$ws.Range("N7").Value = -3229

$ws.Range("H61").Value = 2556.4443
$ws.Range("I61").Value = 2042.1666
$ws.Range("J61").Value = 3585
$ws.Range("K61").Value = 2042.1666
$ws.Range("L61").Value = 3585
$ws.Range("M61").Value = -1840.1666
$ws.Range("N61").Value = -3989

$ws.Range("H113").Value = 2556.4443
$ws.Range("I113").Value = 2042.1666
$ws.Range("J113").Value = 3585
$ws.Range("K113").Value = 2042.1666
$ws.Range("L113").Value = 3585
$ws.Range("M113").Value = 127.8334
$ws.Range("N113").Value = -7925

$ws.Range("H122").Value = 25002142
$ws.Range("I122").Value = 27779822
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 83339466
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -83337016
$ws.Range("N122").Value = -13915

$ws.Range("H126").Value = 2026.3
$ws.Range("I126").Value = 1917.5555
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 5752.666499999999
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -3282.666499999999
$ws.Range("N126").Value = -13955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 933.2222
$ws.Range("I107").Value = 733.3333
$ws.Range("J107").Value = 1333
$ws.Range("K107").Value = 2199.9999
$ws.Range("L107").Value = 3999
$ws.Range("M107").Value = -279.9998999999998
$ws.Range("N107").Value = -7839

$ws.Range("H122").Value = 14446456
$ws.Range("I122").Value = 15296071
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 45888213
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -45885763
$ws.Range("N122").Value = -13900
